$wb = $excel.ActiveWorkbook

# Fix the typo in the first sheet's name: "Timerline" -> "Timeline"
$wb.Worksheets.Item("Timerline").Name = "Timeline"

# Make the Timeline sheet the active/selected sheet (it was previously the
# Quality sheet that was active).
$wb.Worksheets.Item("Timeline").Activate()
